# Fixing some bugs in the plots and report format. Updating example data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# New layout:
#   A: Sample        (letters A..F)
#   B: Endoscopy      (date, yyyy-mm-dd)
#   C: Pathology      (LGD / HGD / NDBE)
#   D: p53 IHC        (0/1)
#   E: GEJ Distance   (number)
# ---------------------------------------------------------------------

# Start from a clean sheet.
$ws.Cells.Clear()

# --- Header row (row 1), bold ------------------------------------------
$ws.Range("A1").Value = "Sample"
$ws.Range("B1").Value = "Endoscopy"
$ws.Range("C1").Value = "Pathology"
$ws.Range("D1").Value = "p53 IHC"
$ws.Range("E1").Value = "GEJ Distance"
$ws.Range("A1:E1").Font.Bold = $true
$ws.Range("B1").NumberFormat = "yyyy\-mm\-dd;@"

# --- Data rows (2-7) -----------------------------------------------------
$data = @(
  @("A", 39096, "LGD",  0, 34),
  @("B", 39255, "LGD",  0, 34),
  @("C", 39255, "LGD",  1, 32),
  @("D", 39390, "HGD",  1, 32),
  @("E", 39390, "NDBE", 0, 30),
  @("F", 39390, "NDBE", 0, 34)
)

$row = 2
foreach ($rec in $data) {
  $ws.Cells.Item($row, 1).Value = $rec[0]
  $ws.Cells.Item($row, 2).Value = $rec[1]
  $ws.Cells.Item($row, 3).Value = $rec[2]
  $ws.Cells.Item($row, 4).Value = $rec[3]
  $ws.Cells.Item($row, 5).Value = $rec[4]
  $row = $row + 1
}

$dataRange = $ws.Range("A2:E7")
$dateRange = $ws.Range("B2:B7")

# Number format for the date column.
$dateRange.NumberFormat = "yyyy\-mm\-dd;@"

# Fonts: column A uses size 12, columns B:E use size 11; both black (not theme).
$ws.Range("A2:A7").Font.Size = 12
$ws.Range("B2:E7").Font.Size = 11
$dataRange.Font.Color = 0

# Center align + unlock all data cells, and give the rows extra height.
$dataRange.HorizontalAlignment = -4108
$dataRange.Locked = $false
$ws.Rows("2:7").RowHeight = 25

# --- View / selection -----------------------------------------------------
$ws.Range("E10").Select()

Write-Output "done"
